$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values need to be swapped between row 13 and row 14
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "X", "Z", "AB")

foreach ($col in $cols) {
    $addr13 = "$col`13"
    $addr14 = "$col`14"
    $val13 = $ws.Range($addr13).Value2
    $val14 = $ws.Range($addr14).Value2
    $ws.Range($addr13).Value = $val14
    $ws.Range($addr14).Value = $val13
}
